# Update FFXIV Leve profit calculations (Brynhildr_Profits) across all job sheets.
# Values refreshed from the latest Universalis market snapshot by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 262.66666
$ws.Range("I2").Value = 170.5
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 170.5
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -57.5
$ws.Range("N2").Value = -1226
$ws.Range("H106").Value = 9374.4
$ws.Range("I106").Value = 9374.4
$ws.Range("K106").Value = 9374.4
$ws.Range("M106").Value = -8743.4
$ws.Range("H115").Value = 562.38464
$ws.Range("I115").Value = 562.38464
$ws.Range("K115").Value = 1687.15392
$ws.Range("M115").Value = -120.15392
$ws.Range("H125").Value = 1066.6666
$ws.Range("I125").Value = 600
$ws.Range("K125").Value = 5400
$ws.Range("M125").Value = -2940
$ws.Range("H137").Value = 26107.666
$ws.Range("I137").Value = 20829
$ws.Range("J137").Value = 36665
$ws.Range("K137").Value = 62487
$ws.Range("L137").Value = 109995
$ws.Range("M137").Value = -59937
$ws.Range("N137").Value = -115095
$ws.Range("H138").Value = 4576.0967
$ws.Range("I138").Value = 9456.833000000001
$ws.Range("J138").Value = 3404.72
$ws.Range("K138").Value = 28370.499
$ws.Range("L138").Value = 10214.16
$ws.Range("M138").Value = -23230.499
$ws.Range("N138").Value = -20494.16

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 300
$ws.Range("K4").Value = 300
$ws.Range("M4").Value = -184
$ws.Range("H74").Value = 5498.5186
$ws.Range("I74").Value = 2573.0444
$ws.Range("K74").Value = 2573.0444
$ws.Range("M74").Value = -1699.0444
$ws.Range("H77").Value = 5498.5186
$ws.Range("I77").Value = 2573.0444
$ws.Range("K77").Value = 12865.222
$ws.Range("M77").Value = -8497.222000000002
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H131").Value = 92500
$ws.Range("I131").Value = 90000
$ws.Range("J131").Value = 95000
$ws.Range("K131").Value = 90000
$ws.Range("L131").Value = 95000
$ws.Range("M131").Value = -84960
$ws.Range("N131").Value = -105080
$ws.Range("H132").Value = 5262.327
$ws.Range("I132").Value = 3276.7273
$ws.Range("J132").Value = 8240.727999999999
$ws.Range("K132").Value = 9830.1819
$ws.Range("L132").Value = 24722.184
$ws.Range("M132").Value = -7300.1819
$ws.Range("N132").Value = -29782.184

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 614.3077
$ws.Range("I80").Value = 993.9
$ws.Range("J80").Value = 377.0625
$ws.Range("K80").Value = 993.9
$ws.Range("L80").Value = 377.0625
$ws.Range("M80").Value = 4.100000000000023
$ws.Range("N80").Value = -2373.0625
$ws.Range("H83").Value = 614.3077
$ws.Range("I83").Value = 993.9
$ws.Range("J83").Value = 377.0625
$ws.Range("K83").Value = 4969.5
$ws.Range("L83").Value = 1885.3125
$ws.Range("M83").Value = 22.5
$ws.Range("N83").Value = -11869.3125
$ws.Range("H94").Value = 5256.2666
$ws.Range("I94").Value = 5284.6665
$ws.Range("J94").Value = 5142.6665
$ws.Range("K94").Value = 5284.6665
$ws.Range("L94").Value = 5142.6665
$ws.Range("M94").Value = -4833.6665
$ws.Range("N94").Value = -6044.6665
$ws.Range("H105").Value = 2074.7036
$ws.Range("I105").Value = 1876.2632
$ws.Range("J105").Value = 2546
$ws.Range("K105").Value = 1876.2632
$ws.Range("L105").Value = 2546
$ws.Range("M105").Value = -129.2632000000001
$ws.Range("N105").Value = -6040
$ws.Range("H134").Value = 3687.25
$ws.Range("I134").Value = 3687.25
$ws.Range("K134").Value = 11061.75
$ws.Range("M134").Value = -8526.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 10290.167
$ws.Range("I122").Value = 1987.75
$ws.Range("K122").Value = 5963.25
$ws.Range("M122").Value = -3513.25
$ws.Range("H132").Value = 3277.8572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 14770.889
$ws.Range("I80").Value = 7002
$ws.Range("K80").Value = 21006
$ws.Range("M80").Value = -20070
$ws.Range("H83").Value = 14770.889
$ws.Range("I83").Value = 7002
$ws.Range("K83").Value = 63018
$ws.Range("M83").Value = -58338
$ws.Range("H94").Value = 11373.25
$ws.Range("J94").Value = 12426.714
$ws.Range("L94").Value = 37280.142
$ws.Range("N94").Value = -38632.142
$ws.Range("H109").Value = 6583.636
$ws.Range("I109").Value = 605
$ws.Range("K109").Value = 1815
$ws.Range("M109").Value = -775
$ws.Range("H113").Value = 855.5
$ws.Range("I113").Value = 599
$ws.Range("J113").Value = 869
$ws.Range("K113").Value = 1797
$ws.Range("L113").Value = 2607
$ws.Range("M113").Value = 373
$ws.Range("N113").Value = -6947
$ws.Range("H119").Value = 14885.571
$ws.Range("I119").Value = 2099.5
$ws.Range("K119").Value = 6298.5
$ws.Range("M119").Value = -1460.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2510.72
$ws.Range("I102").Value = 2960.4443
$ws.Range("J102").Value = 1354.2858
$ws.Range("K102").Value = 2960.4443
$ws.Range("L102").Value = 1354.2858
$ws.Range("M102").Value = -1338.4443
$ws.Range("N102").Value = -4598.2858
$ws.Range("H122").Value = 3248.4
$ws.Range("I122").Value = 3122.7778
$ws.Range("K122").Value = 9368.3334
$ws.Range("M122").Value = -6918.3334
$ws.Range("H132").Value = 12910.25
$ws.Range("I132").Value = 19517.715
$ws.Range("K132").Value = 58553.145
$ws.Range("M132").Value = -56023.145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6200.4287
$ws.Range("I7").Value = 5679.6
$ws.Range("K7").Value = 5679.6
$ws.Range("M7").Value = -5567.6
$ws.Range("H22").Value = 3027.9707
$ws.Range("I22").Value = 2375.9333
$ws.Range("K22").Value = 2375.9333
$ws.Range("M22").Value = -2080.9333
$ws.Range("H27").Value = 3027.9707
$ws.Range("I27").Value = 2375.9333
$ws.Range("K27").Value = 2375.9333
$ws.Range("M27").Value = -2268.9333
$ws.Range("H61").Value = 7163.846
$ws.Range("I61").Value = 7086.3237
$ws.Range("J61").Value = 7691
$ws.Range("K61").Value = 7086.3237
$ws.Range("L61").Value = 7691
$ws.Range("M61").Value = -6884.3237
$ws.Range("N61").Value = -8095
$ws.Range("H113").Value = 7163.846
$ws.Range("I113").Value = 7086.3237
$ws.Range("J113").Value = 7691
$ws.Range("K113").Value = 7086.3237
$ws.Range("L113").Value = 7691
$ws.Range("M113").Value = -4916.3237
$ws.Range("N113").Value = -12031
$ws.Range("H126").Value = 6200.4287
$ws.Range("I126").Value = 5679.6
$ws.Range("K126").Value = 17038.8
$ws.Range("M126").Value = -14568.8
$ws.Range("H136").Value = 2457.0688
$ws.Range("I136").Value = 2093.9583
$ws.Range("K136").Value = 6281.874899999999
$ws.Range("M136").Value = -3731.874899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3133.5557
$ws.Range("I81").Value = 2210.7144
$ws.Range("K81").Value = 4421.4288
$ws.Range("M81").Value = -3360.4288
$ws.Range("H84").Value = 3133.5557
$ws.Range("I84").Value = 2210.7144
$ws.Range("K84").Value = 22107.144
$ws.Range("M84").Value = -16803.144
$ws.Range("H107").Value = 1976.8182
$ws.Range("I107").Value = 1272
$ws.Range("J107").Value = 3210.25
$ws.Range("K107").Value = 3816
$ws.Range("L107").Value = 9630.75
$ws.Range("M107").Value = -1896
$ws.Range("N107").Value = -13470.75
$ws.Range("H113").Value = 27778196
$ws.Range("I113").Value = 666.5
$ws.Range("J113").Value = 55555730
$ws.Range("K113").Value = 1999.5
$ws.Range("L113").Value = 166667190
$ws.Range("M113").Value = 170.5
$ws.Range("N113").Value = -166671530
$ws.Range("H132").Value = 4026.5715
$ws.Range("I132").Value = 3694
$ws.Range("J132").Value = 4359.143
$ws.Range("K132").Value = 11082
$ws.Range("L132").Value = 13077.429
$ws.Range("M132").Value = -8552
$ws.Range("N132").Value = -18137.429
$ws.Range("H136").Value = 1882.625
$ws.Range("I136").Value = 1715
$ws.Range("J136").Value = 2720.75
$ws.Range("K136").Value = 5145
$ws.Range("L136").Value = 8162.25
$ws.Range("M136").Value = -2595
$ws.Range("N136").Value = -13262.25
